$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row 192: 2025-12-05 (serial 45996), station "四方坪站充电量(kw)"
$ws.Range("A192").Value = 45996
$ws.Range("B192").Value = "四方坪站充电量(kw)"

$row192 = @(446.91, 956.67800000000011, 482.76999999999992, 215.10000000000002, 264.23700000000002, 497.9620000000001, 429.95100000000008, 102.524, 113.17399999999999, 180.00400000000002, 66.94, 173.20400000000004, 639.32399999999984, 1731.3229999999996, 719.25599999999986, 492.89300000000003, 367.875, 316.82799999999997, 78.22699999999999, 142.19999999999999, 94.72999999999999, 86.945999999999998, 54.36, 54.599999999999994)

$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $addr = $cols[$i] + "192"
    $ws.Range($addr).Value = $row192[$i]
}

# New data row 193: 2025-12-05 (serial 45996), station "高岭站充电量(kw)"
$ws.Range("A193").Value = 45996
$ws.Range("B193").Value = "高岭站充电量(kw)"

$row193 = @(413.89300000000003, 303.50399999999996, 168.148, 107.506, 66.811999999999998, 87.385999999999996, 296.54499999999996, 171.548, 331.24699999999996, 236.25000000000003, 132.49700000000001, 129.94899999999998, 930.32300000000021, 921.13799999999992, 362.733, 514.00199999999995, 33.347999999999999, 97.497, 14.977, 121.30500000000001, 35.421999999999997, 20.257999999999999, 48.242000000000004, 54.43)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $addr = $cols[$i] + "193"
    $ws.Range($addr).Value = $row193[$i]
}

# Match the author's final selection state
[void]$ws.Range("F198").Select()
